$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5 through 17 (entire rows)
$ws.Range("A5:B17").EntireRow.Delete()

# Update row 2
$ws.Range("A2").Value = "05/13/2021 20:18:13"
$ws.Range("B2").Value = 0.3390000000000001

# Update row 3
$ws.Range("A3").Value = "05/13/2021 20:23:32"
$ws.Range("B3").Value = 0.3279

# Update row 4
$ws.Range("A4").Value = "05/13/2021 20:23:51"
$ws.Range("B4").Value = 3.500999999999999
